# Weekly update for "Comercializadora del Agro de Limarí - Pepino dulce".
# A new week's worth of data (4 quality grades, date 44769) is inserted at
# the top of this supplier's block (rows 433-436). All of the existing
# rows for this supplier (previously rows 433-469) shift down by 4 rows
# (to rows 437-473), which also pushes the sheet's used range/dimension
# from A1:R469 to A1:R473.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the existing block (rows 433-469, columns A-R) into memory
#    BEFORE writing anything, then write it back 4 rows lower. Value2 is
#    read in full first, so the overlapping source/destination ranges do
#    not corrupt each other.
$srcRange = $ws.Range("A433:R469")
$blockValues = $srcRange.Value2

$dstRange = $ws.Range("A437:R473")
$dstRange.Value = $blockValues

# 2) Overwrite rows 433-436 with the new week's data (date 44769).
$newRows = @(
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44769, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Especial", 248, 11000, 12000, 11516, "`$/bandeja 18 kilos", "Provincia de Limarí", 640, 18, "Hortaliza"),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44769, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 400, 9000, 10000, 9500, "`$/bandeja 18 kilos", "Provincia de Limarí", 528, 18, "Hortaliza"),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44769, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 300, 7000, 8000, 7500, "`$/bandeja 18 kilos", "Provincia de Limarí", 417, 18, "Hortaliza"),
    @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44769, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Tercera", 200, 4000, 5000, 4500, "`$/bandeja 18 kilos", "Provincia de Limarí", 250, 18, "Hortaliza")
)

$startRow = 433
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

# 3) The block shift created rows (470-473) that didn't previously exist,
#    so column D there lacks the date number format the rest of the "Fecha"
#    column uses. Re-apply it across the whole touched block (433-473) so
#    every date cell renders/serialises consistently.
$ws.Range("D433:D473").NumberFormat = "YYYY-MM-DD HH:MM:SS"
